$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.958.90'
$ws.Range("E2").Value = '  -0.31%  '

$ws.Range("D3").Value = '1.741.87'
$ws.Range("E3").Value = '  -0.34%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.20%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '249.93'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.63%  '

$ws.Range("E6").Value = '  +0.22%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5137'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.00%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2750'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.37%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06187'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.02%  '

$ws.Range("D10").Value = '1.746.44'
$ws.Range("E10").Value = '  -0.02%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07225'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.13%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.10'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.75%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6496'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.79%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.631'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.45%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '77.59'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.01%  '

$ws.Range("E16").Value = '  +0.17%  '

$ws.Range("E17").Value = '  +0.24%  '

$ws.Range("D18").Value = '25.978.84'
$ws.Range("E18").Value = '  +0.07%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.83'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.79%  '

$ws.Range("E20").Value = '  +0.79%  '

$ws.Range("D21").Value = '1.966.18'
$ws.Range("E21").Value = '  -0.11%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.269'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.14%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.676'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.52%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.384'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.36%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '135.92'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.29%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.508'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.31%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.779'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.77%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '106.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.61%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.951'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.73%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08204'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.41%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.652'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.87%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04701'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.53%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.661'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.76%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9975'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.49%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6243'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.62%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.729'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.62%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01613'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.00%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.913'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.48%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.001'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.25%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '99.98'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.25%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7574'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.28%  '

$ws.Range("E43").Value = '  -1.68%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.021'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.35%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.304'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.58%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1129'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.44%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.53'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.67%  '

$ws.Range("E48").Value = '  -2.15%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.77'
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.514'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.52%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3423'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.36%  '

